# Applies the 31/12/2025 18:27 scrape update to the three schedule sheets:
#   LP1912 (+20 rows), LP1912-215 (+3 rows), 6203-6173 (+1 row).
# Row layout differs per sheet:
#   LP1912:      A(blank) B=Hora_Scrap C=Hora_Llegada D=Linea       E=Minutos F=Parada G=Fecha
#   LP1912-215:  A(blank) B=Fecha       C=Hora_Scrap   D=Hora_Llegada E=Linea  F=Minutos G=Parada
#   6203-6173:   A(blank) B=Fecha       C=Hora_Scrap   D=Hora_Llegada E=Linea  F=Minutos G=Parada
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 31/12/2025 18:27:34"
$ws1.Range("A3").Value = "Total filas: 1220"

$sheet1Data = @(
    @(1202, "18:27:23", "18:31", "16_SANTA ANA", 4, "LP1912", "31/12/2025"),
    @(1203, "18:27:23", "18:33", "23_HERNANDEZ", 6, "LP1912", "31/12/2025"),
    @(1204, "18:27:23", "18:34", "14X44_ABASTO", 7, "LP1912", "31/12/2025"),
    @(1205, "18:27:23", "18:41", "16_P MOR-SANTA ANA", 14, "LP1912", "31/12/2025"),
    @(1206, "18:27:23", "18:44", "14_ABASTO", 17, "LP1912", "31/12/2025"),
    @(1207, "18:27:23", "18:51", "15_ABASTO", 24, "LP1912", "31/12/2025"),
    @(1208, "18:27:23", "18:53", "16_SANTA ANA", 26, "LP1912", "31/12/2025"),
    @(1209, "18:27:23", "19:01", "17_ROMERO", 34, "LP1912", "31/12/2025"),
    @(1210, "18:27:23", "19:03", "23_HERNANDEZ", 36, "LP1912", "31/12/2025"),
    @(1211, "18:27:23", "19:05", "16_SANTA ANA", 38, "LP1912", "31/12/2025"),
    @(1212, "18:27:23", "19:11", "81_EL PELIGRO", 44, "LP1912", "31/12/2025"),
    @(1213, "18:27:23", "19:14", "14_ABASTO", 47, "LP1912", "31/12/2025"),
    @(1214, "18:27:23", "19:21", "215C_EL PATO", 54, "LP1912", "31/12/2025"),
    @(1215, "18:27:23", "19:29", "225_GOMEZ", 62, "LP1912", "31/12/2025"),
    @(1216, "18:27:23", "19:31", "215_EL PELIGRO", 64, "LP1912", "31/12/2025"),
    @(1217, "18:27:23", "19:32", "23_HERNANDEZ", 65, "LP1912", "31/12/2025"),
    @(1218, "18:27:23", "19:44", "11_ETCHEVERRY", 77, "LP1912", "31/12/2025"),
    @(1219, "18:27:23", "19:51", "81_EL PELIGRO", 84, "LP1912", "31/12/2025"),
    @(1220, "18:27:23", "19:59", "14X44_ABASTO", 92, "LP1912", "31/12/2025"),
    @(1221, "18:27:23", "20:01", "215C_EL PATO", 94, "LP1912", "31/12/2025")
)
foreach ($row in $sheet1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 31/12/2025 18:27:34"
$ws2.Range("A3").Value = "Total filas: 83"

$sheet2Data = @(
    @(82, "31/12/2025", "18:27:23", "19:21", "215C_EL PATO", 54, "LP1912"),
    @(83, "31/12/2025", "18:27:23", "19:31", "215_EL PELIGRO", 64, "LP1912"),
    @(84, "31/12/2025", "18:27:23", "20:01", "215C_EL PATO", 94, "LP1912")
)
foreach ($row in $sheet2Data) {
    $r = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
}

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 31/12/2025 18:27:34"
$ws3.Range("A3").Value = "Total filas: 144"

$sheet3Data = ,@(145, "31/12/2025", "18:27:33", "19:11", "215B_LP-P MOR-1 Y 57", 44, "L6173")
foreach ($row in $sheet3Data) {
    $r = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
}

